$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.552.26'
$ws.Range("E2").Value = '  -5.06%  '
$ws.Range("D3").Value = '1.840.75'
$ws.Range("E3").Value = '  -4.35%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").Value = '''312.87'
$ws.Range("E5").Value = '  -3.93%  '
$ws.Range("D6").Value = '''0.9995'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").Value = '''0.4231'
$ws.Range("E7").Value = '  -7.79%  '
$ws.Range("D8").Value = '''0.3629'
$ws.Range("E8").Value = '  -4.99%  '
$ws.Range("D9").Value = '''43.57'
$ws.Range("E9").Value = '  -4.67%  '
$ws.Range("D10").Value = '''0.07198'
$ws.Range("E10").Value = '  -7.14%  '
$ws.Range("D11").Value = '''0.8984'
$ws.Range("E11").Value = '  -8.29%  '
$ws.Range("E12").Value = '  -8.91%  '
$ws.Range("D13").Value = '1.805.45'
$ws.Range("E13").Value = '  -8.42%  '
$ws.Range("D14").Value = '''6.574'
$ws.Range("E14").Value = '  -5.85%  '
$ws.Range("D15").Value = '''5.307'
$ws.Range("E15").Value = '  -7.04%  '
$ws.Range("D16").Value = '''0.06798'
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '''77.00'
$ws.Range("E18").Value = '  -9.27%  '
$ws.Range("D19").Value = '''0.000008893'
$ws.Range("E19").Value = '  -6.45%  '
$ws.Range("D20").Value = '''0.9996'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").Value = '''15.29'
$ws.Range("D22").Value = '27.525.13'
$ws.Range("E22").Value = '  -5.25%  '
$ws.Range("E23").Value = '  -7.94%  '
$ws.Range("D24").Value = '''10.74'
$ws.Range("E24").Value = '  -2.86%  '
$ws.Range("D25").Value = '2.052.10'
$ws.Range("E25").Value = '  -5.63%  '
$ws.Range("D26").Value = '''2.044'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '''151.13'
$ws.Range("E27").Value = '  -4.37%  '
$ws.Range("D28").Value = '''18.21'
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("D29").Value = '''5.264'
$ws.Range("E29").Value = '  -6.42%  '
$ws.Range("D30").Value = '''110.58'
$ws.Range("E30").Value = '  -6.00%  '
$ws.Range("D31").Value = '''1.682'
$ws.Range("E31").Value = '  -8.30%  '
$ws.Range("D32").Value = '''0.08854'
$ws.Range("E32").Value = '  -5.09%  '
$ws.Range("D33").Value = '''0.7723'
$ws.Range("E33").Value = '  -10.45%  '
$ws.Range("D34").Value = '''4.488'
$ws.Range("E34").Value = '  -12.08%  '
$ws.Range("D35").Value = '''2.845'
$ws.Range("E35").Value = '  -5.66%  '
$ws.Range("E36").Value = '  -13.16%  '
$ws.Range("D37").Value = '''1.000'
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").Value = '''0.05381'
$ws.Range("E38").Value = '  -5.64%  '
$ws.Range("D39").Value = '''1.095'
$ws.Range("E39").Value = '  -5.17%  '
$ws.Range("D40").Value = '''0.01919'
$ws.Range("E40").Value = '  -6.56%  '
$ws.Range("D41").Value = '''2.941'
$ws.Range("E41").Value = '  -4.87%  '
$ws.Range("D42").Value = '''0.5033'
$ws.Range("E42").Value = '  -8.70%  '
$ws.Range("D43").Value = '''6.769'
$ws.Range("E43").Value = '  -9.22%  '
$ws.Range("D44").Value = '''0.1633'
$ws.Range("E44").Value = '  -6.87%  '
$ws.Range("D45").Value = '''0.06609'
$ws.Range("E45").Value = '  -4.72%  '
$ws.Range("D46").Value = '''8.203'
$ws.Range("E46").Value = '  -12.21%  '
$ws.Range("D47").Value = '''0.4716'
$ws.Range("E47").Value = '  -9.15%  '
$ws.Range("D48").Value = '''105.41'
$ws.Range("E48").Value = '  -4.96%  '
$ws.Range("D49").Value = '''10.17'
$ws.Range("E49").Value = '  -9.66%  '
$ws.Range("D50").Value = '''0.9990'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").Value = '''1.640'
$ws.Range("E51").Value = '  -6.93%  '
